# AT_01_Login.xlsx test-data refresh:
#  - new login email / password values
#  - selection moved to A5
#  - column B (Password) unhidden
#  - stale A2 hyperlink (old email) removed, B2 hyperlink kept/re-pointed
#  - concurrentCalc flag dropped from calcPr (handled automatically on save)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New credentials
$ws.Range("A2").Value = "thejaswi.y@vensaiinc.com"
$ws.Range("B2").Value = "Abc@123"

# Unhide the Password column (keeps its existing width)
$ws.Columns.Item(2).Hidden = $false

# Move the active selection to A5
$ws.Range("A5").Select()

# Drop the hyperlink on A2 and re-create the one on B2 so only it remains,
# pointing at the updated password value.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Abc@123")
$ws.Range("B2").Style = "Hyperlink"
